$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A38").Value = "Daniele Ruzzenenti "
$ws.Range("B38").Value = "ELIA BATTISTI | U.S. Guarna"
$ws.Range("C38").Value = "Carlo Stedile | MAI UNA GIOIA"
$ws.Range("D38").Value = "Alessandro  Maffei | FC Savignano"
$ws.Range("E38").Value = "Riccardo Zeni | Demobusters"
$ws.Range("F38").Value = "Jacopo Zecchinelli | Vigili del Fusto"
